{"js": "// Update the Courier ID and the Receiver's name / address / phone fields.\n// Every target string is unique in the document, so `body.search(...)`\n// (exact match, no wildcards) followed by `range.insertText(..., \"Replace\")`\n// swaps the text in place while leaving the run/paragraph formatting intact.\n\nasync function replaceText(ctx, findText, newText) {\n  const results = ctx.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await ctx.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await ctx.sync();\n}\n\n// Courier ID (leading spaces before \"RR...\" must be preserved)\nawait replaceText(context, \"RR000000023MA\", \"RR000000019MA\");\n\n// Receiver block\nawait replaceText(context, \"MR.nas ons\", \"MR.qwqw qwqw\");\nawait replaceText(context, \"5dsd\", \"Ain Chegga : erqw\");\nawait replaceText(context, \"5684531\", \"123123\");\n", "ps1": "# Update the Courier ID and the Receiver's name / address / phone fields.\n# Each value in this document is unique, so a plain Find/Replace (no\n# wildcards) targeting the whole document range is safe and keeps all of\n# the surrounding run/paragraph formatting untouched.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\n# Courier ID (leading spaces before \"RR...\" must be preserved)\nReplace-Text \"RR000000023MA\" \"RR000000019MA\"\n\n# Receiver block\nReplace-Text \"MR.nas ons\" \"MR.qwqw qwqw\"\nReplace-Text \"5dsd\" \"Ain Chegga : erqw\"\nReplace-Text \"5684531\" \"123123\"\n"}
